$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03888266680073
$ws.Cells.Item(2, 4).Value = 1.042267144998042
$ws.Cells.Item(2, 5).Value = 1.05607053928482
$ws.Cells.Item(2, 6).Value = 1.062674356131369
$ws.Cells.Item(2, 9).Value = 1.040645688227272
$ws.Cells.Item(2, 10).Value = 1.043977693024201
$ws.Cells.Item(2, 11).Value = 1.045044315519916
$ws.Cells.Item(2, 12).Value = 1.05880925935154
$ws.Cells.Item(2, 13).Value = 1.06539506323659
$ws.Cells.Item(2, 14).Value = 1.01855334139399

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.03984310399646
$ws.Cells.Item(3, 4).Value = 1.043002576049184
$ws.Cells.Item(3, 5).Value = 1.057182019370021
$ws.Cells.Item(3, 6).Value = 1.063817617615216
$ws.Cells.Item(3, 9).Value = 1.040904277526378
$ws.Cells.Item(3, 10).Value = 1.044583051955594
$ws.Cells.Item(3, 11).Value = 1.045590740876615
$ws.Cells.Item(3, 12).Value = 1.059733546959199
$ws.Cells.Item(3, 13).Value = 1.06635236148839
$ws.Cells.Item(3, 14).Value = 1.0187572271746

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.040464664828921
$ws.Cells.Item(4, 4).Value = 1.043478499644876
$ws.Cells.Item(4, 5).Value = 1.057902214138186
$ws.Cells.Item(4, 6).Value = 1.064558262947926
$ws.Cells.Item(4, 9).Value = 1.04107050100455
$ws.Cells.Item(4, 10).Value = 1.044974238528252
$ws.Cells.Item(4, 11).Value = 1.045943697944438
$ws.Cells.Item(4, 12).Value = 1.060332023662407
$ws.Cells.Item(4, 13).Value = 1.066972089172042
$ws.Cells.Item(4, 14).Value = 1.018888906341151

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.040725991007356
$ws.Cells.Item(5, 4).Value = 1.043678589099605
$ws.Cells.Item(5, 5).Value = 1.058205221535299
$ws.Cells.Item(5, 6).Value = 1.064869840431396
$ws.Cells.Item(5, 9).Value = 1.041140117433978
$ws.Cells.Item(5, 10).Value = 1.04513856807497
$ws.Cells.Item(5, 11).Value = 1.046091932871477
$ws.Cells.Item(5, 12).Value = 1.060583719057317
$ws.Cells.Item(5, 13).Value = 1.067232692312183
$ws.Cells.Item(5, 14).Value = 1.018944204502883

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.040769870100311
$ws.Cells.Item(6, 4).Value = 1.043712185628007
$ws.Cells.Item(6, 5).Value = 1.058256111733888
$ws.Cells.Item(6, 6).Value = 1.064922167982235
$ws.Cells.Item(6, 9).Value = 1.041151790853169
$ws.Cells.Item(6, 10).Value = 1.045166152366674
$ws.Cells.Item(6, 11).Value = 1.046116813458901
$ws.Cells.Item(6, 12).Value = 1.060625985434009
$ws.Cells.Item(6, 13).Value = 1.06727645280275
$ws.Cells.Item(6, 14).Value = 1.018953485795874

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.040468156596181
$ws.Cells.Item(7, 4).Value = 1.043481173207208
$ws.Cells.Item(7, 5).Value = 1.057906262004775
$ws.Cells.Item(7, 6).Value = 1.064562425435831
$ws.Cells.Item(7, 9).Value = 1.041071432260108
$ws.Cells.Item(7, 10).Value = 1.044976434800652
$ws.Cells.Item(7, 11).Value = 1.045945679249945
$ws.Cells.Item(7, 12).Value = 1.060335386452406
$ws.Cells.Item(7, 13).Value = 1.066975571090643
$ws.Cells.Item(7, 14).Value = 1.018889645473026

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.039207231373311
$ws.Cells.Item(8, 4).Value = 1.042515676349027
$ws.Cells.Item(8, 5).Value = 1.056445963508266
$ws.Cells.Item(8, 6).Value = 1.06306054465222
$ws.Cells.Item(8, 9).Value = 1.040733307488106
$ws.Cells.Item(8, 10).Value = 1.044182384496406
$ws.Cells.Item(8, 11).Value = 1.045229109859279
$ws.Cells.Item(8, 12).Value = 1.059121543582619
$ws.Cells.Item(8, 13).Value = 1.065718526135091
$ws.Cells.Item(8, 14).Value = 1.018622296863886

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.036986057082567
$ws.Cells.Item(9, 4).Value = 1.040814775002965
$ws.Cells.Item(9, 5).Value = 1.053880342141856
$ws.Cells.Item(9, 6).Value = 1.06042077630492
$ws.Cells.Item(9, 9).Value = 1.040129069370208
$ws.Cells.Item(9, 10).Value = 1.042779198246837
$ws.Cells.Item(9, 11).Value = 1.043961725948965
$ws.Cells.Item(9, 12).Value = 1.056985674335859
$ws.Cells.Item(9, 13).Value = 1.063505695586717
$ws.Cells.Item(9, 14).Value = 1.018149301799634

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.035505793260821
$ws.Cells.Item(10, 4).Value = 1.039681181479719
$ws.Cells.Item(10, 5).Value = 1.05217505253646
$ws.Cells.Item(10, 6).Value = 1.058665468367816
$ws.Cells.Item(10, 9).Value = 1.039720601509367
$ws.Cells.Item(10, 10).Value = 1.041841102516359
$ws.Cells.Item(10, 11).Value = 1.043113678980278
$ws.Cells.Item(10, 12).Value = 1.055563837839129
$ws.Cells.Item(10, 13).Value = 1.062031993591883
$ws.Cells.Item(10, 14).Value = 1.017832713698319

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.034864950005094
$ws.Cells.Item(11, 4).Value = 1.039190414602316
$ws.Cells.Item(11, 5).Value = 1.051437859539488
$ws.Cells.Item(11, 6).Value = 1.057906477693044
$ws.Cells.Item(11, 9).Value = 1.039542395643088
$ws.Cells.Item(11, 10).Value = 1.04143427675317
$ws.Cells.Item(11, 11).Value = 1.042745730916178
$ws.Cells.Item(11, 12).Value = 1.054948660291674
$ws.Cells.Item(11, 13).Value = 1.061394225882357
$ws.Cells.Item(11, 14).Value = 1.017695332054812

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.034626930673876
$ws.Cells.Item(12, 4).Value = 1.039008135764261
$ws.Cells.Item(12, 5).Value = 1.051164214788294
$ws.Cells.Item(12, 6).Value = 1.057624715366021
$ws.Cells.Item(12, 9).Value = 1.039476001481397
$ws.Cells.Item(12, 10).Value = 1.041283070185777
$ws.Cells.Item(12, 11).Value = 1.042608948207443
$ws.Cells.Item(12, 12).Value = 1.054720229113524
$ws.Cells.Item(12, 13).Value = 1.061157383908106
$ws.Cells.Item(12, 14).Value = 1.017644258000268

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.034677985792332
$ws.Cells.Item(13, 4).Value = 1.039047234554702
$ws.Cells.Item(13, 5).Value = 1.051222904301261
$ws.Cells.Item(13, 6).Value = 1.057685147074534
$ws.Cells.Item(13, 9).Value = 1.039490252326746
$ws.Cells.Item(13, 10).Value = 1.04131550872626
$ws.Cells.Item(13, 11).Value = 1.042638293551562
$ws.Cells.Item(13, 12).Value = 1.054769225066542
$ws.Cells.Item(13, 13).Value = 1.06120818490103
$ws.Cells.Item(13, 14).Value = 1.017655215563925

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.034845274876385
$ws.Cells.Item(14, 4).Value = 1.039175347089793
$ws.Cells.Item(14, 5).Value = 1.051415236275579
$ws.Cells.Item(14, 6).Value = 1.057883183855128
$ws.Cells.Item(14, 9).Value = 1.039536911571115
$ws.Cells.Item(14, 10).Value = 1.041421779873427
$ws.Cells.Item(14, 11).Value = 1.042734426664449
$ws.Cells.Item(14, 12).Value = 1.054929776606999
$ws.Cells.Item(14, 13).Value = 1.061374647362439
$ws.Cells.Item(14, 14).Value = 1.017691111161956

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.034948349646273
$ws.Cells.Item(15, 4).Value = 1.03925428330828
$ws.Cells.Item(15, 5).Value = 1.051533762409083
$ws.Cells.Item(15, 6).Value = 1.058005222132198
$ws.Cells.Item(15, 9).Value = 1.039565633301448
$ws.Cells.Item(15, 10).Value = 1.041487244674872
$ws.Cells.Item(15, 11).Value = 1.042793642831456
$ws.Cells.Item(15, 12).Value = 1.05502870740308
$ws.Cells.Item(15, 13).Value = 1.061477217448859
$ws.Cells.Item(15, 14).Value = 1.017713221748138

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.035548326482111
$ws.Cells.Item(16, 4).Value = 1.039713753955719
$ws.Cells.Item(16, 5).Value = 1.052224003069606
$ws.Cells.Item(16, 6).Value = 1.058715862603579
$ws.Cells.Item(16, 9).Value = 1.039732400305785
$ws.Cells.Item(16, 10).Value = 1.041868089087187
$ws.Cells.Item(16, 11).Value = 1.043138082970955
$ws.Cells.Item(16, 12).Value = 1.055604675406572
$ws.Cells.Item(16, 13).Value = 1.062074327614744
$ws.Cells.Item(16, 14).Value = 1.017841825028131

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.035924708600687
$ws.Cells.Item(17, 4).Value = 1.040001991518805
$ws.Cells.Item(17, 5).Value = 1.05265729636477
$ws.Cells.Item(17, 6).Value = 1.059161914834934
$ws.Cells.Item(17, 9).Value = 1.039836651167335
$ws.Cells.Item(17, 10).Value = 1.042106815913799
$ws.Cells.Item(17, 11).Value = 1.043353943891328
$ws.Cells.Item(17, 12).Value = 1.055966095297875
$ws.Cells.Item(17, 13).Value = 1.062448974095561
$ws.Cells.Item(17, 14).Value = 1.017922415127194

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.036144257556243
$ws.Cells.Item(18, 4).Value = 1.040170123913287
$ws.Cells.Item(18, 5).Value = 1.052910145632889
$ws.Cells.Item(18, 6).Value = 1.059422192723799
$ws.Cells.Item(18, 9).Value = 1.039897329953706
$ws.Cells.Item(18, 10).Value = 1.042246000987127
$ws.Cells.Item(18, 11).Value = 1.043479780662939
$ws.Cells.Item(18, 12).Value = 1.056176952398808
$ws.Cells.Item(18, 13).Value = 1.062667533298145
$ws.Cells.Item(18, 14).Value = 1.017969393307156

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.036219120042248
$ws.Cells.Item(19, 4).Value = 1.040227454097788
$ws.Cells.Item(19, 5).Value = 1.052996380481005
$ws.Cells.Item(19, 6).Value = 1.059510958258211
$ws.Cells.Item(19, 9).Value = 1.03991799796774
$ws.Cells.Item(19, 10).Value = 1.042293449269334
$ws.Cells.Item(19, 11).Value = 1.04352267565009
$ws.Cells.Item(19, 12).Value = 1.056248857192664
$ws.Cells.Item(19, 13).Value = 1.062742062146878
$ws.Cells.Item(19, 14).Value = 1.017985406784595

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.035884325156217
$ws.Cells.Item(20, 4).Value = 1.039971065491596
$ws.Cells.Item(20, 5).Value = 1.052610796032472
$ws.Cells.Item(20, 6).Value = 1.059114046980835
$ws.Cells.Item(20, 9).Value = 1.039825479375495
$ws.Cells.Item(20, 10).Value = 1.042081209013127
$ws.Cells.Item(20, 11).Value = 1.04333079142378
$ws.Cells.Item(20, 12).Value = 1.055927313503313
$ws.Cells.Item(20, 13).Value = 1.06240877453105
$ws.Cells.Item(20, 14).Value = 1.017913771535292

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.034796011901119
$ws.Cells.Item(21, 4).Value = 1.03913762075292
$ws.Cells.Item(21, 5).Value = 1.051358594291439
$ws.Cells.Item(21, 6).Value = 1.057824862526012
$ws.Cells.Item(21, 9).Value = 1.03952317712134
$ws.Cells.Item(21, 10).Value = 1.04139048824342
$ws.Cells.Item(21, 11).Value = 1.042706120902882
$ws.Cells.Item(21, 12).Value = 1.054882496160695
$ws.Cells.Item(21, 13).Value = 1.061325626843637
$ws.Cells.Item(21, 14).Value = 1.017680542026838

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034111852671771
$ws.Cells.Item(22, 4).Value = 1.038613681397647
$ws.Cells.Item(22, 5).Value = 1.050572335744303
$ws.Cells.Item(22, 6).Value = 1.057015230280375
$ws.Cells.Item(22, 9).Value = 1.039331947388971
$ws.Cells.Item(22, 10).Value = 1.040955664699796
$ws.Cells.Item(22, 11).Value = 1.042312727165811
$ws.Cells.Item(22, 12).Value = 1.054226001268741
$ws.Cells.Item(22, 13).Value = 1.060644917580378
$ws.Cells.Item(22, 14).Value = 1.017533644625577

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.034474528244891
$ws.Cells.Item(23, 4).Value = 1.03889142366526
$ws.Cells.Item(23, 5).Value = 1.050989046722407
$ws.Cells.Item(23, 6).Value = 1.057444343587439
$ws.Cells.Item(23, 9).Value = 1.039433431791715
$ws.Cells.Item(23, 10).Value = 1.041186223936729
$ws.Cells.Item(23, 11).Value = 1.042521332999012
$ws.Cells.Item(23, 12).Value = 1.054573981518707
$ws.Cells.Item(23, 13).Value = 1.061005745173732
$ws.Cells.Item(23, 14).Value = 1.017611541983927

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.035902572660464
$ws.Cells.Item(24, 4).Value = 1.039985039606492
$ws.Cells.Item(24, 5).Value = 1.052631807170588
$ws.Cells.Item(24, 6).Value = 1.059135676084741
$ws.Cells.Item(24, 9).Value = 1.039830527825849
$ws.Cells.Item(24, 10).Value = 1.042092779855431
$ws.Cells.Item(24, 11).Value = 1.04334125324743
$ws.Cells.Item(24, 12).Value = 1.055944837181614
$ws.Cells.Item(24, 13).Value = 1.06242693887794
$ws.Cells.Item(24, 14).Value = 1.017917677290976

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.037560194123293
$ws.Cells.Item(25, 4).Value = 1.041254442616759
$ws.Cells.Item(25, 5).Value = 1.054542714065638
$ws.Cells.Item(25, 6).Value = 1.061102420495388
$ws.Cells.Item(25, 9).Value = 1.040286275255822
$ws.Cells.Item(25, 10).Value = 1.043142422713255
$ws.Cells.Item(25, 11).Value = 1.044289928227594
$ws.Cells.Item(25, 12).Value = 1.057537482334883
$ws.Cells.Item(25, 13).Value = 1.064077498460743
$ws.Cells.Item(25, 14).Value = 1.018271805160151

